# Updating info.xlsx and template.xlsx | Implementing multiple PM header info
#
# The template previously shipped with one project's hard-coded header data
# (resource phone numbers, cost-center code, TI coordinator name, parcial
# dates) baked directly into the cells. To let the sheet serve multiple
# PMs/projects, those fixed values are cleared out here (leaving the label
# cells/formulas intact) and one leftover label is replaced with a neutral
# "#" placeholder.

$wb = $excel.ActiveWorkbook

$wsSheet = $wb.Worksheets.Item("sheet")
$wsTotal = $wb.Worksheets.Item("TOTALIZADOR")
$wsCab   = $wb.Worksheets.Item("CABEÇALHO")

# --- "sheet" tab: clear the stray phone-number label that was sitting in D1 ---
$wsSheet.Range("D1").Value = $null

# --- TOTALIZADOR tab: clear the stray "COL" label in A3 ---
$wsTotal.Range("A3").Value = $null

# --- CABEÇALHO tab: widen column B to fit the new placeholder text, clear the
#     project-specific parcial dates / cost-center code / TI responsible name,
#     and swap the hard-coded resource contact for a generic "#" ---
$wsCab.Columns.Item(2).ColumnWidth = 22.14
$wsCab.Range("B5").Value = $null
$wsCab.Range("C5").Value = $null
$wsCab.Range("B7").Value = "#"
$wsCab.Range("B10").Value = $null
$wsCab.Range("B12").Value = $null

# --- restore the selection/active-sheet state recorded in each sheet view,
#     finishing on CABEÇALHO so it is the active tab on reopen ---
$wsSheet.Activate() | Out-Null
$wsSheet.Range("E6").Select() | Out-Null

$wsTotal.Activate() | Out-Null
$wsTotal.Range("A3:C6").Select() | Out-Null

$wsCab.Activate() | Out-Null
$wsCab.Range("B5:C5").Select() | Out-Null
